# Trade #47 closed at 2026-02-17 12:48:33 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" roll-up figures to reflect the
# newly closed MarketMaking trade, and appends the trade's own row (#48) to
# both the "All Trades" and "MarketMaking" ledgers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.15   # Current Capital
$summary.Range("B4").Value = 0.14      # Total P&L $
$summary.Range("B5").Value = 0.06      # Total P&L %
$summary.Range("B6").Value = 47        # Total Trades
$summary.Range("B7").Value = 20        # Winning Trades
$summary.Range("B9").Value = 42.55     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.15     # Capital
$status.Range("D4").Value = 47         # Trades
$status.Range("E4").Value = 0.14       # P&L $
$status.Range("F4").Value = 0.15       # P&L %
$status.Range("G4").Value = 42.55      # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append trade #47 as a new row (row 48) on both the "All Trades" and
#    "MarketMaking" ledgers - they carry identical data.
# ---------------------------------------------------------------------------
function Add-TradeRow48($ws) {
    $ws.Range("A48").Value = 47

    # Dates/times that look like dates get auto-converted to serials by
    # Excel; force the cell to Text first so the literal string is kept,
    # then drop the number format back to the sheet default so no stray
    # style is left behind.
    $ws.Range("B48").NumberFormat = "@"
    $ws.Range("B48").Value = "2026-02-17"
    $ws.Range("B48").Style = "Normal"

    $ws.Range("C48").Value = "12:48:27"
    $ws.Range("D48").Value = "MarketMaking"
    $ws.Range("E48").Value = "DOWN"
    $ws.Range("F48").Value = 0.04
    $ws.Range("G48").Value = 0.066984
    $ws.Range("H48").Value = "CLOSED"
    $ws.Range("I48").Value = 67.4593
    $ws.Range("J48").Value = 0.03
    $ws.Range("K48").Value = 100.15
    $ws.Range("L48").Value = 0
    $ws.Range("M48").Value = 0
    $ws.Range("N48").Value = 0.6
    $ws.Range("O48").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P48").Value = "early_exit"
    $ws.Range("Q48").Value = 0.13
}

Add-TradeRow48 $wb.Worksheets.Item("All Trades")
Add-TradeRow48 $wb.Worksheets.Item("MarketMaking")
